# "Update communion columns added" - add Father/Mother/Sponsor columns to
# the first-communion form: Name, Birth Date, Father, Mother, Sponsor 1,
# Sponsor 2, Contact Number, Present Address (A..H), 4 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---- 1. Column widths (A stays; B/C, D/E/F, G, H get new widths) ----
$ws.Columns.Item(2).ColumnWidth = 19.666666666666668
$ws.Columns.Item(3).ColumnWidth = 19.666666666666668
$ws.Columns.Item(4).ColumnWidth = 22
$ws.Columns.Item(5).ColumnWidth = 22
$ws.Columns.Item(6).ColumnWidth = 22
$ws.Columns.Item(7).ColumnWidth = 19
$ws.Columns.Item(8).ColumnWidth = 34

# ---- 2. Banner row (row 1) - extend the centered/blank style to F1:H1 ----
foreach ($col in 6..8) {
    $c = $ws.Cells.Item(1, $col)
    $c.Value = "x"
    $c.HorizontalAlignment = $xlCenter
    $c.ClearContents()
}

# ---- 3. Header row (row 2) ----
$ws.Cells.Item(2, 1).Value = "Name"
$ws.Cells.Item(2, 2).Value = "Birth Date"
$ws.Cells.Item(2, 3).Value = "Father"
$ws.Cells.Item(2, 4).Value = "Mother"
$ws.Cells.Item(2, 5).Value = "Sponsor 1"
$ws.Cells.Item(2, 6).Value = "Sponsor 2"
$ws.Cells.Item(2, 7).Value = "Contact Number"
$ws.Cells.Item(2, 8).Value = "Present Address"

$headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, 8))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = $xlCenter

# trailing blank (but centered) header cells I2:J2
foreach ($col in 9..10) {
    $c = $ws.Cells.Item(2, $col)
    $c.Value = "x"
    $c.HorizontalAlignment = $xlCenter
    $c.ClearContents()
}

# ---- 4. Data rows 3-6 ----
$names = @("John Mark Victorino", "Shekinah Joy Victorino", "Jedidiah Victorino", "Kyla Faith Victorino")
$father = "Mark Lucas Victorino"
$mother = "Sarah Lee"
$sponsor1 = "Richard Bowers"
$sponsor2 = "Mikee Corbito"
$address = "Naic, Cavite"
$birthDate = 34899
$contactNumber = 9062268483

for ($i = 0; $i -lt 4; $i++) {
    $row = 3 + $i

    $ws.Cells.Item($row, 1).Value = $names[$i]

    $ws.Cells.Item($row, 2).Value = $birthDate
    $ws.Cells.Item($row, 2).NumberFormat = "yyyy\-mm\-dd;@"

    $ws.Cells.Item($row, 3).Value = $father
    $ws.Cells.Item($row, 4).Value = $mother
    $ws.Cells.Item($row, 5).Value = $sponsor1
    $ws.Cells.Item($row, 6).Value = $sponsor2

    $ws.Cells.Item($row, 7).Value = $contactNumber
    $ws.Cells.Item($row, 7).NumberFormat = "00000000000"

    $ws.Cells.Item($row, 8).Value = $address
}

# ---- 5. Merge the banner row across the new width ----
$ws.Range("A1:H1").Merge()

# ---- 6. Selection shown when the file was last saved ----
$ws.Range("E8").Select()
